# "fixes for working hrs constraint"
#
# The Quantity column (D) on the orders sheet was storing raw "working
# hours" values of 35 for every order. The fix rescales those to the
# correct units: most rows become 3500, while a few rows (the ones whose
# constraint was off by a further factor of ten) become 35000.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowValues = @{
    2  = 3500
    3  = 3500
    4  = 3500
    5  = 3500
    6  = 3500
    7  = 3500
    8  = 3500
    9  = 3500
    10 = 3500
    11 = 35000
    12 = 3500
    13 = 3500
    14 = 35000
    15 = 3500
    16 = 35000
    17 = 3500
    18 = 3500
    19 = 3500
    20 = 3500
    21 = 35000
    22 = 3500
    23 = 3500
    24 = 35000
    25 = 3500
    26 = 3500
    27 = 3500
    28 = 3500
}

foreach ($row in $rowValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $rowValues[$row]
}

# Scroll the view down a bit and move the active selection, matching
# where the user was working when they made the fix.
$window = $excel.ActiveWindow
$window.ScrollRow = 7
$window.ScrollColumn = 1

$ws.Range("D11").Select()
